$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics per row (matches canonical OOXML diff)
# Row 2
$ws.Range("G2").Value = 44.29505033333334
$ws.Range("H2").Value = 132.885151
$ws.Range("I2").Value = 0.9830698162761968
$ws.Range("J2").Value = 0.9830698162761969
$ws.Range("M2").Value = 46.63275166666667
$ws.Range("N2").Value = 139.898255
$ws.Range("O2").Value = 0.9158911059585902
$ws.Range("P2").Value = 0.9158911059585902
$ws.Range("Q2").Value = 2065.600082256834
$ws.Range("R2").Value = 18590.40074031151
$ws.Range("S2").Value = 0.9003849012637141
$ws.Range("T2").Value = 0.9003849012637142
# Row 3
$ws.Range("G3").Value = 44.29505033333334
$ws.Range("H3").Value = 132.885151
$ws.Range("I3").Value = 0.9830698162761968
$ws.Range("J3").Value = 0.9830698162761969
$ws.Range("O3").Value = 0.05441917700612491
$ws.Range("P3").Value = 0.05441917700612491
$ws.Range("Q3").Value = 122.7310274866706
$ws.Range("R3").Value = 1104.579247380035
$ws.Range("S3").Value = 0.05349785034131305
$ws.Range("T3").Value = 0.05349785034131306
# Row 4
$ws.Range("G4").Value = 44.29505033333334
$ws.Range("H4").Value = 132.885151
$ws.Range("I4").Value = 0.9830698162761968
$ws.Range("J4").Value = 0.9830698162761969
$ws.Range("M4").Value = 0.849605
$ws.Range("N4").Value = 2.548815
$ws.Range("O4").Value = 0.01668667696558362
$ws.Range("P4").Value = 0.01668667696558362
$ws.Range("Q4").Value = 37.63329623845167
$ws.Range("R4").Value = 338.6996661460651
$ws.Range("S4").Value = 0.01640416845881654
$ws.Range("T4").Value = 0.01640416845881654
# Row 5
$ws.Range("G5").Value = 44.29505033333334
$ws.Range("H5").Value = 132.885151
$ws.Range("I5").Value = 0.9830698162761968
$ws.Range("J5").Value = 0.9830698162761969
$ws.Range("M5").Value = 0.662052
$ws.Range("N5").Value = 1.986156
$ws.Range("O5").Value = 0.01300304006970129
$ws.Range("P5").Value = 0.0130030400697013
$ws.Range("Q5").Value = 29.325626663284
$ws.Range("R5").Value = 263.930639969556
$ws.Range("S5").Value = 0.01278289621235328
$ws.Range("T5").Value = 0.01278289621235328
# Row 6
$ws.Range("I6").Value = 0.006814145293655052
$ws.Range("J6").Value = 0.006814145293655053
$ws.Range("M6").Value = 46.63275166666667
$ws.Range("N6").Value = 139.898255
$ws.Range("O6").Value = 0.9158911059585902
$ws.Range("P6").Value = 0.9158911059585902
$ws.Range("Q6").Value = 14.31770037696833
$ws.Range("R6").Value = 128.859303392715
$ws.Range("S6").Value = 0.006241015069168248
$ws.Range("T6").Value = 0.006241015069168249
# Row 7
$ws.Range("I7").Value = 0.006814145293655052
$ws.Range("J7").Value = 0.006814145293655053
$ws.Range("O7").Value = 0.05441917700612491
$ws.Range("P7").Value = 0.05441917700612491
$ws.Range("S7").Value = 0.0003708201788808673
$ws.Range("T7").Value = 0.0003708201788808674
# Row 8
$ws.Range("I8").Value = 0.006814145293655052
$ws.Range("J8").Value = 0.006814145293655053
$ws.Range("M8").Value = 0.849605
$ws.Range("N8").Value = 2.548815
$ws.Range("O8").Value = 0.01668667696558362
$ws.Range("P8").Value = 0.01668667696558362
$ws.Range("Q8").Value = 0.260855072755
$ws.Range("R8").Value = 2.347695654795
$ws.Range("S8").Value = 0.0001137054413117738
$ws.Range("T8").Value = 0.0001137054413117738
# Row 9
$ws.Range("I9").Value = 0.006814145293655052
$ws.Range("J9").Value = 0.006814145293655053
$ws.Range("M9").Value = 0.662052
$ws.Range("N9").Value = 1.986156
$ws.Range("O9").Value = 0.01300304006970129
$ws.Range("P9").Value = 0.0130030400697013
$ws.Range("Q9").Value = 0.203270487612
$ws.Range("R9").Value = 1.829434388508
$ws.Range("S9").Value = 0.00008860460429416314
$ws.Range("T9").Value = 0.00008860460429416315
# Row 10
$ws.Range("G10").Value = 0.3685326666666667
$ws.Range("H10").Value = 1.105598
$ws.Range("I10").Value = 0.008179093108268589
$ws.Range("J10").Value = 0.008179093108268589
$ws.Range("M10").Value = 46.63275166666667
$ws.Range("N10").Value = 139.898255
$ws.Range("O10").Value = 0.9158911059585902
$ws.Range("P10").Value = 0.9158911059585902
$ws.Range("Q10").Value = 17.18569232572111
$ws.Range("R10").Value = 154.67123093149
$ws.Range("S10").Value = 0.007491158632670401
$ws.Range("T10").Value = 0.007491158632670401
# Row 11
$ws.Range("G11").Value = 0.3685326666666667
$ws.Range("H11").Value = 1.105598
$ws.Range("I11").Value = 0.008179093108268589
$ws.Range("J11").Value = 0.008179093108268589
$ws.Range("O11").Value = 0.05441917700612491
$ws.Range("P11").Value = 0.05441917700612491
$ws.Range("Q11").Value = 1.021116185714444
$ws.Range("R11").Value = 9.19004567143
$ws.Range("S11").Value = 0.0004450995156084447
$ws.Range("T11").Value = 0.0004450995156084447
# Row 12
$ws.Range("G12").Value = 0.3685326666666667
$ws.Range("H12").Value = 1.105598
$ws.Range("I12").Value = 0.008179093108268589
$ws.Range("J12").Value = 0.008179093108268589
$ws.Range("M12").Value = 0.849605
$ws.Range("N12").Value = 2.548815
$ws.Range("O12").Value = 0.01668667696558362
$ws.Range("P12").Value = 0.01668667696558362
$ws.Range("Q12").Value = 0.3131071962633334
$ws.Range("R12").Value = 2.817964766370001
$ws.Range("S12").Value = 0.0001364818845691092
$ws.Range("T12").Value = 0.0001364818845691092
# Row 13
$ws.Range("G13").Value = 0.3685326666666667
$ws.Range("H13").Value = 1.105598
$ws.Range("I13").Value = 0.008179093108268589
$ws.Range("J13").Value = 0.008179093108268589
$ws.Range("M13").Value = 0.662052
$ws.Range("N13").Value = 1.986156
$ws.Range("O13").Value = 0.01300304006970129
$ws.Range("P13").Value = 0.0130030400697013
$ws.Range("Q13").Value = 0.243987789032
$ws.Range("R13").Value = 2.195890101288
$ws.Range("S13").Value = 0.0001063530754206342
$ws.Range("T13").Value = 0.0001063530754206342
# Row 14
$ws.Range("E14").Value = 2.0
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08727466666666667
$ws.Range("H14").Value = 0.261824
$ws.Range("I14").Value = 0.001936945321879485
$ws.Range("J14").Value = 0.001936945321879485
$ws.Range("M14").Value = 46.63275166666667
$ws.Range("N14").Value = 139.898255
$ws.Range("O14").Value = 0.9158911059585902
$ws.Range("P14").Value = 0.9158911059585902
$ws.Range("Q14").Value = 4.069857857457778
$ws.Range("R14").Value = 36.62872071712
$ws.Range("S14").Value = 0.001774030993037519
$ws.Range("T14").Value = 0.001774030993037519
# Row 15
$ws.Range("E15").Value = 2.0
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08727466666666667
$ws.Range("H15").Value = 0.261824
$ws.Range("I15").Value = 0.001936945321879485
$ws.Range("J15").Value = 0.001936945321879485
$ws.Range("O15").Value = 0.05441917700612491
$ws.Range("P15").Value = 0.05441917700612491
$ws.Range("Q15").Value = 0.2418173008711111
$ws.Range("R15").Value = 2.17635570784
$ws.Range("S15").Value = 0.0001054069703225453
$ws.Range("T15").Value = 0.0001054069703225453
# Row 16
$ws.Range("E16").Value = 2.0
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08727466666666667
$ws.Range("H16").Value = 0.261824
$ws.Range("I16").Value = 0.001936945321879485
$ws.Range("J16").Value = 0.001936945321879485
$ws.Range("M16").Value = 0.849605
$ws.Range("N16").Value = 2.548815
$ws.Range("O16").Value = 0.01668667696558362
$ws.Range("P16").Value = 0.01668667696558362
$ws.Range("Q16").Value = 0.07414899317333334
$ws.Range("R16").Value = 0.6673409385600001
$ws.Range("S16").Value = 0.00003232118088620136
$ws.Range("T16").Value = 0.00003232118088620136
# Row 17
$ws.Range("E17").Value = 2.0
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08727466666666667
$ws.Range("H17").Value = 0.261824
$ws.Range("I17").Value = 0.001936945321879485
$ws.Range("J17").Value = 0.001936945321879485
$ws.Range("M17").Value = 0.662052
$ws.Range("N17").Value = 1.986156
$ws.Range("O17").Value = 0.01300304006970129
$ws.Range("P17").Value = 0.0130030400697013
$ws.Range("Q17").Value = 0.057780367616
$ws.Range("R17").Value = 0.520023308544
$ws.Range("S17").Value = 0.00002518617763321941
$ws.Range("T17").Value = 0.00002518617763321942
